# Apply updated cryptocurrency price/volume data to match the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Formula = "'39.888.08"
$ws.Cells.Item(2, 5).Value = '  -0.52%  '
$ws.Cells.Item(3, 4).Formula = "'2.210.08"
$ws.Cells.Item(3, 5).Value = '  -1.28%  '
$ws.Cells.Item(4, 5).Value = '  -0.14%  '
$ws.Cells.Item(5, 4).Formula = "'293.29"
$ws.Cells.Item(5, 5).Value = '  -0.78%  '
$ws.Cells.Item(6, 4).Formula = "'87.17"
$ws.Cells.Item(6, 5).Value = '  +0.27%  '
$ws.Cells.Item(7, 4).Formula = "'0.508"
$ws.Cells.Item(7, 5).Value = '  -1.65%  '
$ws.Cells.Item(8, 5).Value = '  -0.19%  '
$ws.Cells.Item(9, 4).Formula = "'0.473"
$ws.Cells.Item(9, 5).Value = '  -0.27%  '
$ws.Cells.Item(10, 4).Formula = "'0.0775"
$ws.Cells.Item(10, 5).Value = '  -2.58%  '
$ws.Cells.Item(11, 4).Formula = "'29.74"
$ws.Cells.Item(11, 5).Value = '  -5.25%  '
$ws.Cells.Item(12, 4).Formula = "'48.91"
$ws.Cells.Item(12, 5).Value = '  +3.61%  '
$ws.Cells.Item(13, 4).Formula = "'0.111"
$ws.Cells.Item(13, 5).Value = '  +2.85%  '
$ws.Cells.Item(14, 4).Formula = "'6.43"
$ws.Cells.Item(14, 5).Value = '  -1.07%  '
$ws.Cells.Item(15, 4).Formula = "'2.553.70"
$ws.Cells.Item(15, 5).Value = '  -1.43%  '
$ws.Cells.Item(16, 4).Formula = "'13.69"
$ws.Cells.Item(16, 5).Value = '  -3.47%  '
$ws.Cells.Item(17, 4).Formula = "'2.200.63"
$ws.Cells.Item(17, 5).Value = '  -1.63%  '
$ws.Cells.Item(18, 4).Formula = "'0.727"
$ws.Cells.Item(18, 5).Value = '  -0.65%  '
$ws.Cells.Item(19, 4).Formula = "'39.821.92"
$ws.Cells.Item(19, 5).Value = '  -0.56%  '
$ws.Cells.Item(20, 4).Formula = "'0.0₃0882"
$ws.Cells.Item(20, 5).Value = '  -1.15%  '
$ws.Cells.Item(21, 4).Formula = "'11.17"
$ws.Cells.Item(21, 5).Value = '  +1.96%  '
$ws.Cells.Item(22, 4).Formula = "'5.76"
$ws.Cells.Item(22, 5).Value = '  -1.06%  '
$ws.Cells.Item(23, 4).Formula = "'65.19"
$ws.Cells.Item(23, 5).Value = '  -0.65%  '
$ws.Cells.Item(24, 4).Formula = "'235.57"
$ws.Cells.Item(24, 5).Value = '  -0.08%  '
$ws.Cells.Item(25, 5).Value = '  +0.07%  '
$ws.Cells.Item(26, 4).Formula = "'2.44"
$ws.Cells.Item(26, 5).Value = '  -1.28%  '
$ws.Cells.Item(27, 5).Value = '  -2.82%  '
$ws.Cells.Item(28, 4).Formula = "'22.47"
$ws.Cells.Item(28, 5).Value = '  -1.81%  '
$ws.Cells.Item(29, 4).Formula = "'2.06"
$ws.Cells.Item(29, 5).Value = '  -7.82%  '
$ws.Cells.Item(30, 4).Formula = "'9.16"
$ws.Cells.Item(30, 5).Value = '  -0.88%  '
$ws.Cells.Item(31, 4).Formula = "'155.25"
$ws.Cells.Item(31, 5).Value = '  +1.80%  '
$ws.Cells.Item(32, 4).Formula = "'31.68"
$ws.Cells.Item(32, 5).Value = '  -5.42%  '
$ws.Cells.Item(33, 4).Formula = "'1.00"
$ws.Cells.Item(33, 5).Value = '  -0.13%  '
$ws.Cells.Item(34, 4).Formula = "'4.87"
$ws.Cells.Item(34, 5).Value = '  -0.42%  '
$ws.Cells.Item(35, 4).Formula = "'0.0709"
$ws.Cells.Item(35, 5).Value = '  -1.50%  '
$ws.Cells.Item(36, 4).Formula = "'2.35"
$ws.Cells.Item(36, 5).Value = '  -1.07%  '
$ws.Cells.Item(37, 4).Formula = "'2.83"
$ws.Cells.Item(37, 5).Value = '  +4.05%  '
$ws.Cells.Item(38, 4).Formula = "'0.111"
$ws.Cells.Item(38, 5).Value = '  -0.44%  '
$ws.Cells.Item(39, 4).Formula = "'15.40"
$ws.Cells.Item(39, 5).Value = '  -6.80%  '
$ws.Cells.Item(40, 5).Value = '  -3.88%  '
$ws.Cells.Item(41, 4).Formula = "'1.65"
$ws.Cells.Item(41, 5).Value = '  -2.87%  '
$ws.Cells.Item(42, 4).Formula = "'2.124.28"
$ws.Cells.Item(42, 5).Value = '  +4.19%  '
$ws.Cells.Item(43, 4).Formula = "'3.74"
$ws.Cells.Item(43, 5).Value = '  -3.03%  '
$ws.Cells.Item(44, 5).Value = '  -7.11%  '
$ws.Cells.Item(45, 2).Value = 'VeChain'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(45, 4).Formula = "'0.0265"
$ws.Cells.Item(45, 5).Value = '  -2.68%  '
$ws.Cells.Item(46, 2).Value = 'EnergySwap'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(46, 4).Formula = "'17.48"
$ws.Cells.Item(46, 5).Value = '  +7.28%  '
$ws.Cells.Item(47, 4).Formula = "'9.69"
$ws.Cells.Item(47, 5).Value = '  -3.40%  '
$ws.Cells.Item(48, 4).Formula = "'2.64"
$ws.Cells.Item(48, 5).Value = '  +2.59%  '
$ws.Cells.Item(49, 4).Formula = "'2.421.74"
$ws.Cells.Item(49, 5).Value = '  -1.56%  '
$ws.Cells.Item(50, 4).Formula = "'1.45"
$ws.Cells.Item(50, 5).Value = '  -1.02%  '
$ws.Cells.Item(51, 2).Value = 'Aave'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(51, 4).Formula = "'87.97"
$ws.Cells.Item(51, 5).Value = '  -1.60%  '
